# Applies the "Updated cryptos list" refresh described in the commit diff:
# - most rows keep their coin but get refreshed Price (D) / Volume(1h) (E) figures
# - four row-pairs (5/6, 28/29, 42/43, 50/51) had their coins swapped in ranking order,
#   so coin name / link / price / volume are rewritten for both rows of each pair
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "43.949.61"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.354.28"
$ws.Range("E3").Value = "  -0.46%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 (BNB)
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'240.75"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6 (XRP)
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "'0.668"
$ws.Range("E6").Value = "  -3.60%  "

# Row 7 (Solana)
$ws.Range("D7").Value = "'73.67"
$ws.Range("E7").Value = "  -0.73%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.603"
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  -2.64%  "

# Row 11 (OKB)
$ws.Range("D11").Value = "'58.78"
$ws.Range("E11").Value = "  +1.52%  "

# Row 12 (Avalanche)
$ws.Range("D12").Value = "'34.24"
$ws.Range("E12").Value = "  +7.88%  "

# Row 13 (Polkadot)
$ws.Range("D13").Value = "'7.36"
$ws.Range("E13").Value = "  -2.20%  "

# Row 14 (TRON)
$ws.Range("E14").Value = "  +0.00%  "

# Row 15 (WrappedliquidstakedEther2.0)
$ws.Range("D15").Value = "2.705.77"
$ws.Range("E15").Value = "  -0.47%  "

# Row 16 (Chainlink)
$ws.Range("D16").Value = "'16.44"
$ws.Range("E16").Value = "  -3.63%  "

# Row 17 (Polygon)
$ws.Range("D17").Value = "'0.915"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18 (WrappedEther)
$ws.Range("D18").Value = "2.355.32"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19 (WrappedBTC)
$ws.Range("D19").Value = "43.851.59"
$ws.Range("E19").Value = "  -0.77%  "

# Row 20 (ShibaInu)
$ws.Range("D20").Value = "'0.0000102"
$ws.Range("E20").Value = "  -1.69%  "

# Row 21 (Uniswap)
$ws.Range("E21").Value = "  -0.11%  "

# Row 22 (Litecoin)
$ws.Range("D22").Value = "'77.75"
$ws.Range("E22").Value = "  -1.23%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").Value = "'257.30"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24 (ImmutableX)
$ws.Range("D24").Value = "'1.92"
$ws.Range("E24").Value = "  +15.95%  "

# Row 25 (Dai)
$ws.Range("E25").Value = "  -0.04%  "

# Row 26 (WEMIXToken)
$ws.Range("D26").Value = "'3.75"
$ws.Range("E26").Value = "  -0.20%  "

# Row 27 (PancakeSwap)
$ws.Range("D27").Value = "'2.51"
$ws.Range("E27").Value = "  -2.22%  "

# Row 28 (Toncoin)
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.38"
$ws.Range("E28").Value = "  +2.77%  "

# Row 29 (Cosmos)
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'10.63"
$ws.Range("E29").Value = "  -1.76%  "

# Row 30 (EthereumClassic)
$ws.Range("E30").Value = "  -0.20%  "

# Row 31 (Monero)
$ws.Range("D31").Value = "'177.94"
$ws.Range("E31").Value = "  +1.73%  "

# Row 32 (Kaspa)
$ws.Range("E32").Value = "  -0.53%  "

# Row 33 (Stellar)
$ws.Range("E33").Value = "  +0.28%  "

# Row 34 (Hedera)
$ws.Range("E34").Value = "  -0.23%  "

# Row 35 (Filecoin)
$ws.Range("E35").Value = "  -3.79%  "

# Row 36 (InternetComputer(DFINITY))
$ws.Range("D36").Value = "'5.46"
$ws.Range("E36").Value = "  +0.60%  "

# Row 37 (RenderToken)
$ws.Range("D37").Value = "'3.81"
$ws.Range("E37").Value = "  -2.73%  "

# Row 38 (LidoDAOToken)
$ws.Range("E38").Value = "  -3.26%  "

# Row 39 (THORChain)
$ws.Range("D39").Value = "'6.40"
$ws.Range("E39").Value = "  -2.41%  "

# Row 40 (VeChain)
$ws.Range("E40").Value = "  +1.22%  "

# Row 41 (MultiversX)
$ws.Range("D41").Value = "'67.06"
$ws.Range("E41").Value = "  +25.76%  "

# Row 42 (FTXToken)
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "'5.20"
$ws.Range("E42").Value = "  +15.85%  "

# Row 43 (Cronos)
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.111"
$ws.Range("E43").Value = "  +10.65%  "

# Row 44 (FraxShare)
$ws.Range("D44").Value = "'9.13"
$ws.Range("E44").Value = "  +0.23%  "

# Row 45 (InjectiveProtocol)
$ws.Range("D45").Value = "'19.10"
$ws.Range("E45").Value = "  -0.49%  "

# Row 46 (Algorand)
$ws.Range("D46").Value = "'0.203"
$ws.Range("E46").Value = "  +1.77%  "

# Row 47 (NEARProtocol)
$ws.Range("D47").Value = "'2.53"
$ws.Range("E47").Value = "  +0.89%  "

# Row 48 (TrustWalletToken)
$ws.Range("E48").Value = "  -0.11%  "

# Row 49 (BinanceUSD)
$ws.Range("E49").Value = "  -0.01%  "

# Row 50 (ARBITRUM)
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'1.16"
$ws.Range("E50").Value = "  -1.82%  "

# Row 51 (Aave)
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'99.46"
$ws.Range("E51").Value = "  -1.88%  "
